$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-11 Sunday", "2024-02-12 Monday"),
    @("786÷6=131, 0", "303÷3=101, 0"),
    @("843÷6=140, 3", "748÷5=149, 3"),
    @("617÷2=308, 1", "223÷2=111, 1"),
    @("566÷6=94, 2", "637÷7=91, 0"),
    @("449÷8=56, 1", "690÷9=76, 6"),
    @("463÷3=154, 1", "929÷5=185, 4"),
    @("970÷3=323, 1", "466÷8=58, 2"),
    @("861÷9=95, 6", "517÷3=172, 1"),
    @("417÷3=139, 0", "116÷9=12, 8"),
    @("622÷8=77, 6", "399÷7=57, 0"),
    @("935÷7=133, 4", "985÷4=246, 1"),
    @("178÷5=35, 3", "468÷6=78, 0"),
    @("939÷6=156, 3", "213÷9=23, 6"),
    @("239÷3=79, 2", "795÷2=397, 1"),
    @("739÷7=105, 4", "433÷3=144, 1"),
    @("172÷7=24, 4", "711÷4=177, 3"),
    @("161÷2=80, 1", "488÷5=97, 3"),
    @("421÷9=46, 7", "144÷5=28, 4"),
    @("519÷2=259, 1", "537÷6=89, 3"),
    @("787÷8=98, 3", "509÷3=169, 2"),
    @("146÷2=73, 0", "900÷6=150, 0"),
    @("679÷2=339, 1", "688÷3=229, 1"),
    @("998÷6=166, 2", "495÷5=99, 0"),
    @("488÷8=61, 0", "465÷3=155, 0"),
    @("869÷9=96, 5", "303÷2=151, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
